$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update standings table (rows 5-18) with new stats after the extra round of games ---

$ws.Range("C5").Value = "ISsoft"
$ws.Range("D5").Value = 11
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = "725 - 586"
$ws.Range("H5").Value = 20

$ws.Range("C6").Value = "Эра-Недвижимости плюс"
$ws.Range("D6").Value = 11
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = "818 - 655"
$ws.Range("H6").Value = 20

$ws.Range("C7").Value = "БГУФК"
$ws.Range("D7").Value = 11
$ws.Range("E7").Value = 8
$ws.Range("F7").Value = 3
$ws.Range("G7").Value = "756 - 603"
$ws.Range("H7").Value = 19

$ws.Range("C8").Value = "GOLDEN HILL"
$ws.Range("D8").Value = 11
$ws.Range("E8").Value = 8
$ws.Range("F8").Value = 3
$ws.Range("G8").Value = "786 - 730"
$ws.Range("H8").Value = 19

$ws.Range("C9").Value = "Грушвиль"
$ws.Range("D9").Value = 11
$ws.Range("E9").Value = 8
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = "905 - 738"
$ws.Range("H9").Value = 19

$ws.Range("C10").Value = "ОПЛАТИ"
$ws.Range("D10").Value = 11
$ws.Range("E10").Value = 8
$ws.Range("F10").Value = 3
$ws.Range("G10").Value = "810 - 688"
$ws.Range("H10").Value = 19

$ws.Range("C11").Value = "Mapogo males"
$ws.Range("D11").Value = 11
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = "831 - 794"
$ws.Range("H11").Value = 18

$ws.Range("C12").Value = "SIRIUS"
$ws.Range("D12").Value = 11
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 6
$ws.Range("G12").Value = "731 - 643"
$ws.Range("H12").Value = 16

$ws.Range("C13").Value = "Стрела"
$ws.Range("D13").Value = 11
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 7
$ws.Range("G13").Value = "671 - 705"
$ws.Range("H13").Value = 15

$ws.Range("C14").Value = "VSS"
$ws.Range("D14").Value = 11
$ws.Range("E14").Value = 4
$ws.Range("F14").Value = 7
$ws.Range("G14").Value = "689 - 746"
$ws.Range("H14").Value = 15

$ws.Range("C15").Value = "Eagles"
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 8
$ws.Range("G15").Value = "636 - 687"
$ws.Range("H15").Value = 14

$ws.Range("C16").Value = "NORD"
$ws.Range("D16").Value = 11
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = "572 - 853"
$ws.Range("H16").Value = 14

$ws.Range("C17").Value = "ЛФК"
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = "613 - 794"
$ws.Range("H17").Value = 12

$ws.Range("C18").Value = "Минск 7х"
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = "487 - 808"
$ws.Range("H18").Value = 11

# --- Append two new match-day blocks at the bottom of the results log (rows 111-119) ---

# Merge the new row ranges first (while cells are still blank/default) so that copying
# the bordered formatting afterwards does not get redistributed across the merge.
$ws.Range("B111:H111").Merge()
$ws.Range("B112:H112").Merge()
$ws.Range("B113:H113").Merge()
$ws.Range("B114:H114").Merge()
$ws.Range("B115:H115").Merge()
$ws.Range("B116:H116").Merge()
$ws.Range("B117:H117").Merge()
$ws.Range("B118:H118").Merge()
$ws.Range("B119:H119").Merge()

# Copy formatting from the existing "date" row (row 102) onto the new date rows
$ws.Range("B102:H102").Copy()
$ws.Range("B111:H111").PasteSpecial(-4122)
$ws.Range("B102:H102").Copy()
$ws.Range("B115:H115").PasteSpecial(-4122)

# Copy formatting from an existing "match result" row (row 103) onto the new match rows
$ws.Range("B103:H103").Copy()
$ws.Range("B112:H112").PasteSpecial(-4122)
$ws.Range("B103:H103").Copy()
$ws.Range("B113:H113").PasteSpecial(-4122)
$ws.Range("B103:H103").Copy()
$ws.Range("B114:H114").PasteSpecial(-4122)
$ws.Range("B103:H103").Copy()
$ws.Range("B116:H116").PasteSpecial(-4122)
$ws.Range("B103:H103").Copy()
$ws.Range("B117:H117").PasteSpecial(-4122)
$ws.Range("B103:H103").Copy()
$ws.Range("B118:H118").PasteSpecial(-4122)
$ws.Range("B103:H103").Copy()
$ws.Range("B119:H119").PasteSpecial(-4122)

# Row heights for the match result rows
$ws.Range("B112:H112").RowHeight = 19.95
$ws.Range("B113:H113").RowHeight = 19.95
$ws.Range("B114:H114").RowHeight = 19.95
$ws.Range("B116:H116").RowHeight = 19.95
$ws.Range("B117:H117").RowHeight = 19.95
$ws.Range("B118:H118").RowHeight = 19.95
$ws.Range("B119:H119").RowHeight = 19.95

# Values: new match day for 2025-03-01
$ws.Range("B111").Value = 45717
$ws.Range("B112").Value = "Mapogo males - Грушвиль 81:100 (16:30, БНТУ)"
$ws.Range("B113").Value = "Минск 7х - Стрела 50:71 (18:00, БНТУ)"
$ws.Range("B114").Value = "ОПЛАТИ - ISsoft 54:72 (19:30, БНТУ)"

# Values: new match day for 2025-03-02
$ws.Range("B115").Value = 45718
$ws.Range("B116").Value = "NORD - ЛФК 64:62 (11:00, БНТУ)"
$ws.Range("B117").Value = "SIRIUS - Эра-Недвижимости плюс 54:55 (12:30, БНТУ)"
$ws.Range("B118").Value = "GOLDEN HILL - БГУФК 58:62 (14:00, БНТУ)"
$ws.Range("B119").Value = "Eagles - VSS 58:68 (15:30, БНТУ)"
